$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Change B59 to a true numeric value (was stored as text "2")
$ws.Range("B59").Value = 2

# Add new row 60 with data
$ws.Range("A60").Value = "Ying Tang"
$ws.Range("B60").Value = "'3"
$ws.Range("C60").Value = "We will add"
$ws.Range("D60").Value = "FBK"
$ws.Range("E60").Value = "WRI"
$ws.Range("F60").Value = "f6da2ad4-28ad-4a7e-bf94-2041c47bfd2f"
$ws.Range("G60").Value = "rk07ZXZRb_annotated.xlsx"
$ws.Range("H60").Value = "We will add this to the discussion to the paper."
